# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The underlying "K" (strikeouts, column G) values are recomputed from the
# source pitching-log calculation and rewritten here as literal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$kValues = @{
    2  = 5
    3  = 3
    4  = 2
    5  = 6
    6  = 2
    7  = 2
    8  = 1
    9  = 4
    10 = 3
    11 = 0
    12 = 2
    13 = 3
    14 = 2
    15 = 5
    16 = 2
    17 = 3
    18 = 3
    19 = 6
    20 = 6
    21 = 6
    22 = 6
    23 = 10
    24 = 6
    25 = 7
    26 = 6
    27 = 2
    28 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
